# "UPDATE TEST CASE FILE"
#
# The tester (HET PANKITKUMAR PARIKH) opened the workbook, switched to the
# "manipulating" sheet and filled in their name in the merged "Tester name: "
# header cell (A1), leaving that sheet as the active tab on save.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("manipulating")

# Fill in the tester's name in the merged header cell (A1:F1).
$ws.Range("A1").Value = "Tester name:  HET PANKITKUMAR PARIKH"

# Make "manipulating" the active/selected tab (was "tokenizing").
$ws.Activate()
